$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Aclass")

# 1. Insert a new header row at the top (old row1 -> row2, old row2 -> row3)
$ws.Rows.Item(1).Insert()

# 2. Insert a new column B (a date column) - everything from old B.. shifts to C..
$ws.Columns.Item(2).Insert()

# 3. New header row with country names
$ws.Cells.Item(1,1).Value = "India"
$ws.Cells.Item(1,2).Value = "Canada"
$ws.Cells.Item(1,3).Value = "Japan"

# 4. Row 2 (was row 1): add new date value in B2, fix up D2/J2/L2
$ws.Cells.Item(2,2).Value = 44247.84171296297
$ws.Cells.Item(2,2).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(2,4).Value = "ENUM_VAL1_NOT_THE_SAME"
$ws.Cells.Item(2,10).Value = 4
$ws.Cells.Item(2,12).Value = "1h3m0.001s"

# 5. Row 3 (was row 2): add new date value in B3, fix up L3
$ws.Cells.Item(3,2).Value = -693593
$ws.Cells.Item(3,2).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item(3,12).Value = "0s"
